$wb = $excel.ActiveWorkbook

# ---- Sheet "Ruteo": fill in the router IP addresses (VTP Noroeste config) ----
$ws1 = $wb.Worksheets.Item("Ruteo")

$ws1.Range("C6").Value  = "130.45.1.1 "
$ws1.Range("C7").Value  = "130.45.43.1"
$ws1.Range("C8").Value  = "130.45.85.1"
$ws1.Range("C9").Value  = "130.45.128.1"
$ws1.Range("C10").Value = "130.45.172.1"
$ws1.Range("C11").Value = "130.45.225.1"

# C11 used to carry its own distinct bottom border; align it with the rest
# of the IP column (C6:C10) which only has a thin right border.
$ws1.Range("C11").Borders.Item(9).LineStyle = -4142

# Widen the IP column so the longer values fit comfortably.
$ws1.Columns.Item(3).ColumnWidth = 48

# Update the on-screen selection left after the edit.
$ws1.Range("C7").Select()

# ---- Sheet "VLAN": restore it as the active tab/selection ----
$ws2 = $wb.Worksheets.Item("VLAN")
$ws2.Range("A3:D3").Select()
$ws2.Activate()
